$d = $word.ActiveDocument

# Delete paragraphs that contain an inline picture, except for the very
# first inline picture (the "featured" image at the top of the document).
# Walk the InlineShapes collection from the end so that deleting a
# paragraph does not invalidate the indices of the remaining shapes.
for ($i = $d.InlineShapes.Count; $i -ge 2; $i--) {
    $shp = $d.InlineShapes.Item($i)
    $para = $shp.Range.Paragraphs.Item(1)
    $rng = $para.Range
    $rng.MoveEnd(1, 1)
    $rng.Delete()
}
